$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title row
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 12:18"

# Update country data rows that changed between refreshes
$ws.Range("D4").Value = 2903381
$ws.Range("E4").Value = 2453802
$ws.Range("B14").Value = 343203
$ws.Range("C14").Value = 2133
$ws.Range("D14").Value = 297486
$ws.Range("E14").Value = 26078
$ws.Range("G14").Value = 147
$ws.Range("H14").Value = 19639
$ws.Range("A42").Value = "Rumania"
$ws.Range("B42").Value = 70461
$ws.Range("C42").Value = 1087
$ws.Range("D42").Value = 32587
$ws.Range("E42").Value = 34883
$ws.Range("G42").Value = 37
$ws.Range("H42").Value = 2991
$ws.Range("A43").Value = "Bielorrusia"
$ws.Range("B43").Value = 69424
$ws.Range("D43").Value = 66747
$ws.Range("E43").Value = 2070
$ws.Range("H43").Value = 607
$ws.Range("B60").Value = 37596
$ws.Range("C60").Value = 45
$ws.Range("E60").Value = 9055
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 1375
$ws.Range("B87").Value = 9676
$ws.Range("C87").Value = 38
$ws.Range("D87").Value = 8705
$ws.Range("E87").Value = 731
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 240
$ws.Range("B98").Value = 7731
$ws.Range("C98").Value = 11
$ws.Range("E98").Value = 348
$ws.Range("A128").Value = "Lituania"
$ws.Range("C128").Value = 30
$ws.Range("D128").Value = 1704
$ws.Range("E128").Value = 631
$ws.Range("H128").Value = 81
$ws.Range("A129").Value = "Eslovenia"
$ws.Range("B129").Value = 2416
$ws.Range("C129").Value = 15
$ws.Range("D129").Value = 2051
$ws.Range("E129").Value = 236
$ws.Range("H129").Value = 129
$ws.Range("B191").Value = 138
$ws.Range("C191").Value = 5
$ws.Range("E191").Value = 36
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
